$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a brand-new, empty paragraph before the first (bullet) paragraph.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphBefore()

$newFirst = $d.Paragraphs.Item(1)
$newFirstXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newFirst.Range.InsertXML($newFirstXml)

# Force an explicit zero indentation (left=0, firstLine=0) to survive serialization.
$newFirst2 = $d.Paragraphs.Item(1)
$newFirst2.Format.LeftIndent = 0
$newFirst2.Format.FirstLineIndent = 0

# ---------------------------------------------------------------------------
# 2. Update the (now second) paragraph: "In Read_AIA_094_data.ipynb..."
#    -> "Get all files within a certain period..."  font size 16 -> 24.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2Xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="720" w:hanging="360"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Get all files within a certain period into a single pandas dataframe. For example, the entire month of January 2015. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($p2Xml)

# ---------------------------------------------------------------------------
# 3. Update the (now third / last) paragraph: "We need to also make a
#    module..." -> "Immediate goal - is to create a df..." ; ilvl 1 -> 0;
#    indent 1440 -> 720; drop the trailing empty run.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3Range = $p3.Range
$p3Xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="720" w:hanging="360"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Immediate goal - is to create a df with start time end time and flare class, which can be used to tag the AIA images as flare or not flare.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3Range.InsertXML($p3Xml)

# InsertXML on the very last paragraph of the body leaves behind a stray
# empty paragraph (Word always needs a terminating paragraph mark); merge it
# back into the paragraph we just wrote so the body ends with exactly one
# paragraph mark again.
$count = $d.Paragraphs.Count
if ($count -gt 3) {
    $secondLast = $d.Paragraphs.Item($count - 1)
    $lastPara = $d.Paragraphs.Item($count)
    $mergeRange = $d.Range($secondLast.Range.End - 1, $lastPara.Range.End)
    $mergeRange.Delete()
}
